# Sync attendance_reports: swap the order of "dnasr281@gmail.com, System"
# to "System, dnasr281@gmail.com" in column G ("Recorded By") wherever it
# appears with both names combined.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$oldValue = "dnasr281@gmail.com, System"
$newValue = "System, dnasr281@gmail.com"

$searchRange = $ws.Range("G1:G235")

$firstFound = $searchRange.Find($oldValue)
if ($firstFound) {
    $firstAddress = $firstFound.Address()
    $current = $firstFound
    $addresses = @()

    while ($true) {
        $addresses += $current.Address()
        $current = $searchRange.FindNext($current)
        if ((-not $current) -or ($current.Address() -eq $firstAddress)) {
            break
        }
    }

    foreach ($addr in $addresses) {
        $ws.Range($addr).Value = $newValue
    }
}
